$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price values formatted as plain text (e.g. "306.24"), so
# force text format to prevent Excel from auto-converting them to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Cells.Item(2, 4).Value = "47.093.34"
$ws.Cells.Item(2, 5).Value = "  +5.38%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "2.339.59"
$ws.Cells.Item(3, 5).Value = "  +4.41%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  -0.74%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "306.24"
$ws.Cells.Item(5, 5).Value = "  +0.22%  "

# Row 6
$ws.Cells.Item(6, 4).Value = "97.03"
$ws.Cells.Item(6, 5).Value = "  +4.49%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  +1.73%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  -0.55%  "

# Row 9
$ws.Cells.Item(9, 5).Value = "  +5.15%  "

# Row 10
$ws.Cells.Item(10, 4).Value = "35.80"
$ws.Cells.Item(10, 5).Value = "  +3.56%  "

# Row 11
$ws.Cells.Item(11, 5).Value = "  +1.59%  "

# Row 12
$ws.Cells.Item(12, 5).Value = "  +4.58%  "

# Row 13
$ws.Cells.Item(13, 5).Value = "  -0.49%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "2.695.85"
$ws.Cells.Item(14, 5).Value = "  +4.42%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "2.339.68"
$ws.Cells.Item(15, 5).Value = "  +4.28%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "14.15"
$ws.Cells.Item(16, 5).Value = "  +5.26%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "46.914.34"
$ws.Cells.Item(18, 5).Value = "  +5.59%  "

# Row 19
$ws.Cells.Item(19, 4).Value = "13.58"
$ws.Cells.Item(19, 5).Value = "  +16.93%  "

# Row 20
$ws.Cells.Item(20, 5).Value = "  +2.38%  "

# Row 21
$ws.Cells.Item(21, 4).Value = "6.17"
$ws.Cells.Item(21, 5).Value = "  +0.53%  "

# Row 22
$ws.Cells.Item(22, 5).Value = "  +3.55%  "

# Row 23
$ws.Cells.Item(23, 4).Value = "250.97"
$ws.Cells.Item(23, 5).Value = "  +5.71%  "

# Row 24
$ws.Cells.Item(24, 4).Value = "2.99"
$ws.Cells.Item(24, 5).Value = "  +2.34%  "

# Row 25
$ws.Cells.Item(25, 5).Value = "  +1.76%  "

# Row 26
$ws.Cells.Item(26, 5).Value = "  -0.40%  "

# Row 27
$ws.Cells.Item(27, 4).Value = "42.74"
$ws.Cells.Item(27, 5).Value = "  +15.58%  "

# Row 28
$ws.Cells.Item(28, 4).Value = "2.25"
$ws.Cells.Item(28, 5).Value = "  -1.40%  "

# Row 29
$ws.Cells.Item(29, 4).Value = "9.88"

# Row 30
$ws.Cells.Item(30, 4).Value = "20.24"
$ws.Cells.Item(30, 5).Value = "  +1.82%  "

# Row 31
$ws.Cells.Item(31, 5).Value = "  +0.39%  "

# Row 32
$ws.Cells.Item(32, 4).Value = "0.0818"
$ws.Cells.Item(32, 5).Value = "  +5.69%  "

# Row 33
$ws.Cells.Item(33, 4).Value = "147.82"
$ws.Cells.Item(33, 5).Value = "  -0.53%  "

# Row 34
$ws.Cells.Item(34, 4).Value = "2.62"
$ws.Cells.Item(34, 5).Value = "  +0.08%  "

# Row 35
$ws.Cells.Item(35, 4).Value = "3.17"
$ws.Cells.Item(35, 5).Value = "  +1.10%  "

# Row 36
$ws.Cells.Item(36, 5).Value = "  +5.10%  "

# Row 37
$ws.Cells.Item(37, 5).Value = "  +1.28%  "

# Row 38
$ws.Cells.Item(38, 5).Value = "  -1.83%  "

# Row 39
$ws.Cells.Item(39, 5).Value = "  +7.31%  "

# Row 40
$ws.Cells.Item(40, 5).Value = "  +6.36%  "

# Row 41
$ws.Cells.Item(41, 5).Value = "  +2.41%  "

# Row 42
$ws.Cells.Item(42, 4).Value = "13.95"
$ws.Cells.Item(42, 5).Value = "  -6.30%  "

# Row 43
$ws.Cells.Item(43, 5).Value = "  -0.82%  "

# Row 44
$ws.Cells.Item(44, 4).Value = "1.98"
$ws.Cells.Item(44, 5).Value = "  +13.55%  "

# Row 45
$ws.Cells.Item(45, 4).Value = "1.827.96"
$ws.Cells.Item(45, 5).Value = "  +1.36%  "

# Row 46
$ws.Cells.Item(46, 4).Value = "89.22"
$ws.Cells.Item(46, 5).Value = "  +10.13%  "

# Row 47
$ws.Cells.Item(47, 4).Value = "75.02"
$ws.Cells.Item(47, 5).Value = "  +9.74%  "

# Row 48
$ws.Cells.Item(48, 5).Value = "  +5.22%  "

# Row 49
$ws.Cells.Item(49, 4).Value = "99.00"
$ws.Cells.Item(49, 5).Value = "  +1.72%  "

# Row 50 (THORChain -> MultiversX)
$ws.Cells.Item(50, 2).Value = "MultiversX"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Cells.Item(50, 4).Value = "55.25"
$ws.Cells.Item(50, 5).Value = "  +3.33%  "

# Row 51 (MultiversX -> THORChain)
$ws.Cells.Item(51, 2).Value = "THORChain"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Cells.Item(51, 4).Value = "4.85"
$ws.Cells.Item(51, 5).Value = "  +0.88%  "
